$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 11 (the "krigelron@gmail.com" review row). This shifts the
#    old row 12 ("hermanliran@gmail.com" review) up to become row 11, and
#    shrinks the used range from A1:F12 to A1:F11.
$ws.Rows(11).Delete()

# 2. Add the new review's recovery email in D4 ("nachushayinc@gmail.com"),
#    matching the look of the other hyperlink-styled email cells (Arial 10,
#    blue) without disturbing the shared style table used elsewhere.
$ws.Range("D4").Value = "nachushayinc@gmail.com"
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 10
$ws.Range("D4").Font.Color = 16711680

# 3. Rebuild the hyperlinks collection to match the new layout. (The engine's
#    Hyperlinks collection does not automatically re-target itself when rows
#    are deleted/shifted, so we clear it and re-add every link in the
#    correct, final order.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:leviadlevi22@gmail.com", "", "", "leviadlevi22@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:hermanliran@gmail.com", "", "", "hermanliran@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:freelancernachus@gmail.com", "", "", "freelancernachus@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:nachushayinc@gmail.com", "", "", "nachushayinc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:nevilgreen@gmail.com", "", "", "nevilgreen@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:stclerari834@gmail.com", "", "", "stclerari834@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:stcydouel274@gmail.com", "", "", "stcydouel274@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:sinuspai@gmail.com", "", "", "sinuspai@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:hermanliran@gmail.com", "", "", "hermanliran@gmail.com")

# 4. Update the selected cell shown in the sheet view to A11 (matches the
#    author re-selecting the now-last row after the edit).
[void]$ws.Range("A11").Select()
